# Append 10 new order lines (rows 12-21) to the "Orders" sheet and update
# the "Summary" sheet's G2 tracking-number cell.
#
# Each new row only populates column C (FlowerName) and, for all but the
# last row, column F (Number). Values are written as text (matching the
# existing "numberStoredAsText" cells in this sheet) by forcing the
# NumberFormat to "@" before the assignment, then resetting the style back
# to Normal afterwards so no stray per-cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$newRows = @(
  @{ Name = "221_朱丽叶塔_Julieta_Rosa rugosa Thunb._10stems";                               Number = "10" },
  @{ Name = "244_繁星_undefined_Rosa rugosa Thunb._10stems";                                 Number = "10" },
  @{ Name = "433_红豆_Hypericum red_undefined_1bunch";                                       Number = "2"  },
  @{ Name = "434_绿豆_Hypericum green_undefined_1bunch";                                     Number = "2"  },
  @{ Name = "527_白豆_Hypericum white_undefined_1bunch";                                     Number = "2"  },
  @{ Name = "548_白星花_tweedia white_undefined_1bunch";                                     Number = "20" },
  @{ Name = "413_风铃花淡紫色_Canterbury Bells`nlight purple_undefined_1bunch";              Number = "10" },
  @{ Name = "412_紫罗兰粉_violet pink_undefined_1bunch";                                     Number = "5"  },
  @{ Name = "256_奇迹女神_Miracle Goddess_Rosa rugosa Thunb._20stems";                       Number = "10" },
  @{ Name = "203_佛罗伊德_Floyd_Rosa rugosa Thunb._20stems";                                 Number = $null }
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $entry = $newRows[$i]

    Set-TextValue $ws.Range("C$r") $entry.Name

    if ($null -ne $entry.Number) {
        Set-TextValue $ws.Range("F$r") $entry.Number
    }
}

# Summary sheet: the shipment tracking-number string in G2 grew longer.
$wsSummary = $wb.Worksheets.Item("Summary")
Set-TextValue $wsSummary.Range("G2") "0202035101010101010101022220105100"
